# "base code for using DEAP" — CM Freshmen timetable restructuring.
# Rebuilds Wed/Thu (row 6), removes a Wed Sociology slot (rows 14-19),
# and reshuffles the 14:00 and 16:00 blocks (rows 26-39), incl. merged
# cell housekeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Stash the three donor fills/formats we need later into a scratch
#    row (100) far below the real data, so later writes don't depend
#    on write-order among the cells being edited. Deleted in step 8.
#      A100 <- green  (92D050) donor, taken from D26 (style untouched by this edit)
#      B100 <- yellow (FFFF00) donor, taken from E34 (style untouched by this edit)
#      C100 <- white  (FFFFFF) donor, taken from F6  (style untouched by this edit)
# ---------------------------------------------------------------------
$ws.Range("D26").Copy($ws.Range("A100"))
$ws.Range("E34").Copy($ws.Range("B100"))
$ws.Range("F6").Copy($ws.Range("C100"))

# ---------------------------------------------------------------------
# 2) Row 6 (09:00 block): shift English Writing out, add a second
#    Russian (Intermediate) section, and swap Russian Beginner -> Sociology.
# ---------------------------------------------------------------------
$ws.Range("C6").Value = "Russian Language (Elementary Level)`n09:00-10:30`nroom:Creative room: 104"
$ws.Range("D6").Value = "Russian Language (Intermediate Level)`n09:00-10:30`nroom:Creative room: 104"
$ws.Range("E6").Value = "Sociology`n09:00-10:30`nroom:Creative room: 104"
# F6 is unchanged (still "Russian Language (Intermediate Level)" 09:00-10:30).

# ---------------------------------------------------------------------
# 3) Row 14-19 (11:00 block): drop the Wednesday Sociology column
#    entirely (cell + the blank placeholder cells below it).
# ---------------------------------------------------------------------
foreach ($r in 14..19) {
    $c = $ws.Range("D$r")
    $c.Value = ""
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# 4) Row 26-31 (14:00 block): Mathematics II (Monday) is dropped from
#    this slot; Wednesday now carries Mathematics II (new room), and
#    Thursday/Friday both carry English Writing for Media (new room).
# ---------------------------------------------------------------------
foreach ($r in 26..31) {
    $c = $ws.Range("B$r")
    $c.Value = ""
    $c.Style = "Normal"
}

$ws.Range("C100").Copy($ws.Range("D26"))
$ws.Range("D26").Value = "Mathematics II`n14:00-15:30`nroom:White classroom: 202"

$ws.Range("B100").Copy($ws.Range("E26"))
$ws.Range("E26").Value = "English Writing for Media`n14:00-15:30`nroom:Yellow classroom: 209"

$ws.Range("B100").Copy($ws.Range("F26"))
$ws.Range("F26").Value = "English Writing for Media`n14:00-15:30`nroom:Yellow classroom: 209"

# ---------------------------------------------------------------------
# 5) Row 34-39 (16:00 block): Monday gains Russian Beginner, Tuesday
#    keeps Russian Beginner (new room-consistent copy), Thursday becomes
#    Media Literacy, Friday becomes Russian Elementary.
# ---------------------------------------------------------------------
$ws.Range("B100").Copy($ws.Range("B34"))
$ws.Range("B34").Value = "Russian Language (Beginner Level)`n16:00-17:30`nroom:Yellow classroom: 209"

$ws.Range("B100").Copy($ws.Range("C34"))
$ws.Range("C34").Value = "Russian Language (Beginner Level)`n16:00-17:30`nroom:Yellow classroom: 209"

$ws.Range("A100").Copy($ws.Range("E34"))
$ws.Range("E34").Value = "Media Literacy`n16:00-17:30`nroom:Green classroom: 204"

$ws.Range("B100").Copy($ws.Range("F34"))
$ws.Range("F34").Value = "Russian Language (Elementary Level)`n16:00-17:30`nroom:Yellow classroom: 209"

# New blank placeholder cells under the new Monday 16:00 block (B35:B39),
# matching the styling used by the analogous C/E/F columns.
foreach ($r in 35..38) {
    $ws.Range("C$r").Copy($ws.Range("B$r"))
    $ws.Range("B$r").Value = ""
}
$ws.Range("C39").Copy($ws.Range("B39"))
$ws.Range("B39").Value = ""

# ---------------------------------------------------------------------
# 6) Merged cells: drop B26:B31 and D14:D19, add B34:B39. `.Copy()`
#    above can silently break pre-existing merges on its target cells,
#    so re-assert the *entire* final merge set unconditionally.
# ---------------------------------------------------------------------
$ws.Range("B3:B4").Merge()
$ws.Range("B1:G1").Merge()
$ws.Range("F3:F4").Merge()
$ws.Range("C6:C11").Merge()
$ws.Range("F34:F39").Merge()
$ws.Range("D6:D11").Merge()
$ws.Range("C14:C19").Merge()
$ws.Range("F6:F11").Merge()
$ws.Range("A3:A4").Merge()
$ws.Range("D3:D4").Merge()
$ws.Range("E34:E39").Merge()
$ws.Range("B14:B19").Merge()
$ws.Range("F26:F31").Merge()
$ws.Range("E6:E11").Merge()
$ws.Range("B34:B39").Merge()
$ws.Range("E26:E31").Merge()
$ws.Range("C3:C4").Merge()
$ws.Range("E3:E4").Merge()
$ws.Range("D26:D31").Merge()
$ws.Range("C34:C39").Merge()

# ---------------------------------------------------------------------
# 7) Drop the scratch row entirely so it leaves no trace in the saved
#    worksheet (dimension, empty <c> remnants, etc.).
# ---------------------------------------------------------------------
$ws.Rows.Item(100).Delete()
